$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(456).Insert()

$ws.Cells.Item(456, 1).Value = 3
$ws.Cells.Item(456, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(456, 3).Value = "Coquimbo"
$ws.Cells.Item(456, 4).Value = 44855
$ws.Cells.Item(456, 5).Value = 5
$ws.Cells.Item(456, 6).Value = 100112037
$ws.Cells.Item(456, 7).Value = "Cebollín"
$ws.Cells.Item(456, 8).Value = "Sin especificar"
$ws.Cells.Item(456, 9).Value = "Primera"
$ws.Cells.Item(456, 10).Value = 120
$ws.Cells.Item(456, 11).Value = 3500
$ws.Cells.Item(456, 12).Value = 3500
$ws.Cells.Item(456, 13).Value = 3500
$ws.Cells.Item(456, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(456, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(456, 16).Value = 97
$ws.Cells.Item(456, 17).Value = 36
$ws.Cells.Item(456, 18).Value = "Hortaliza"
